$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column C (Glucose [g/L]) for rows 2-20
$ws.Range("C2:C20").ClearContents()

# Clear column J (added Glucose [g]) for rows 12-20
$ws.Range("J12:J20").ClearContents()

# Clear column K (Consumed_Glucose [g]) for rows 2-21
$ws.Range("K2:K21").ClearContents()

# Clear column L (Consumed_Glucose [g/L]) for rows 2-21
$ws.Range("L2:L21").ClearContents()

# Clear column O (Yxs [gx/gs]) for rows 3-21
$ws.Range("O3:O21").ClearContents()

# Clear column P (qs 1 - ds/dt/x - [1/h]) for rows 3-21
$ws.Range("P3:P21").ClearContents()
